$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the 14 new columns at their final positions (left-to-right so each
# insert only affects columns at/after the current target, matching how the
# final layout lines up with the untouched pre-existing columns).
$ws.Range("B1").EntireColumn.Insert()
$ws.Range("E1").EntireColumn.Insert()
$ws.Range("F1").EntireColumn.Insert()
$ws.Range("G1").EntireColumn.Insert()
$ws.Range("I1").EntireColumn.Insert()
$ws.Range("J1").EntireColumn.Insert()
$ws.Range("K1").EntireColumn.Insert()
$ws.Range("L1").EntireColumn.Insert()
$ws.Range("N1").EntireColumn.Insert()
$ws.Range("P1").EntireColumn.Insert()
$ws.Range("R1").EntireColumn.Insert()
$ws.Range("T1").EntireColumn.Insert()
$ws.Range("U1").EntireColumn.Insert()
$ws.Range("W1").EntireColumn.Insert()

# The insert copies formatting from the column to the left, so strip that
# back off the data rows (2-6) of each newly inserted column; row 1 keeps
# the bold/bordered header style already present on the sheet.
$ws.Range("B2:B6").ClearFormats()
$ws.Range("E2:E6").ClearFormats()
$ws.Range("F2:F6").ClearFormats()
$ws.Range("G2:G6").ClearFormats()
$ws.Range("I2:I6").ClearFormats()
$ws.Range("J2:J6").ClearFormats()
$ws.Range("K2:K6").ClearFormats()
$ws.Range("L2:L6").ClearFormats()
$ws.Range("N2:N6").ClearFormats()
$ws.Range("P2:P6").ClearFormats()
$ws.Range("R2:R6").ClearFormats()
$ws.Range("T2:T6").ClearFormats()
$ws.Range("U2:U6").ClearFormats()
$ws.Range("W2:W6").ClearFormats()

# Fill header + data for each newly inserted column (cells left unset stay
# blank, matching the blanks already present for some existing trainees).
$ws.Range("B1").Value = "Alexis Rainey"
$ws.Range("B2").Value = 3.6123
$ws.Range("B3").Value = 2.8192
$ws.Range("B4").Value = 3.0696
$ws.Range("B5").Value = 2.8389
$ws.Range("B6").Value = 3.7651

$ws.Range("E1").Value = "Curley"
$ws.Range("E3").Value = 1.3708
$ws.Range("E4").Value = 1.2718
$ws.Range("E5").Value = 1.1801
$ws.Range("E6").Value = 1.3523

$ws.Range("F1").Value = "Doyle"
$ws.Range("F2").Value = 3.9429
$ws.Range("F3").Value = 2.2349
$ws.Range("F4").Value = 3.3665
$ws.Range("F5").Value = 2.7004
$ws.Range("F6").Value = 3.4779

$ws.Range("G1").Value = "Espona"
$ws.Range("G2").Value = 3.9843
$ws.Range("G3").Value = 2.131
$ws.Range("G4").Value = 2.7494
$ws.Range("G5").Value = 2.8337
$ws.Range("G6").Value = 3.5205

$ws.Range("I1").Value = "Hackman"
$ws.Range("I2").Value = 3.4173
$ws.Range("I3").Value = 1.964
$ws.Range("I4").Value = 3.1906
$ws.Range("I5").Value = 2.2144
$ws.Range("I6").Value = 3.4541

$ws.Range("J1").Value = "Holzman"
$ws.Range("J2").Value = 3.2803
$ws.Range("J3").Value = 2.0907
$ws.Range("J4").Value = 3.138
$ws.Range("J5").Value = 2.0565
$ws.Range("J6").Value = 2.3033

$ws.Range("K1").Value = "Hughes"
$ws.Range("K3").Value = 1.4038
$ws.Range("K4").Value = 2.3211
$ws.Range("K5").Value = 1.7743
$ws.Range("K6").Value = 1.7384

$ws.Range("L1").Value = "Johnson"
$ws.Range("L2").Value = 3.4826
$ws.Range("L3").Value = 2.6531
$ws.Range("L4").Value = 3.6441
$ws.Range("L5").Value = 2.6041
$ws.Range("L6").Value = 3.0425

$ws.Range("N1").Value = "McFadden"
$ws.Range("N2").Value = 3.8007
$ws.Range("N3").Value = 2.1425
$ws.Range("N4").Value = 2.9674
$ws.Range("N5").Value = 2.851
$ws.Range("N6").Value = 3.0195

$ws.Range("P1").Value = "Myers"
$ws.Range("P2").Value = 3.5482
$ws.Range("P3").Value = 1.7411
$ws.Range("P4").Value = 3.0397

$ws.Range("R1").Value = "Reilly"
$ws.Range("R3").Value = 2.1894
$ws.Range("R4").Value = 2.917
$ws.Range("R5").Value = 2.2343
$ws.Range("R6").Value = 2.7926

$ws.Range("T1").Value = "Streib"
$ws.Range("T5").Value = 1.0637
$ws.Range("T6").Value = 1.3416

$ws.Range("U1").Value = "Tollaksen"
$ws.Range("U2").Value = 3.8899
$ws.Range("U3").Value = 2.1552
$ws.Range("U4").Value = 3.6476
$ws.Range("U5").Value = 2.8995

$ws.Range("W1").Value = "Yanovich"
$ws.Range("W2").Value = 3.1199
$ws.Range("W3").Value = 2.0553
